$d = $word.ActiveDocument

$replacements = @(
    @("789÷9=", "426÷6="),
    @("946÷8=", "601÷4="),
    @("316÷4=", "649÷5="),
    @("416÷8=", "401÷7="),
    @("719÷9=", "398÷3="),
    @("203÷7=", "617÷2="),
    @("855÷9=", "895÷3="),
    @("657÷9=", "657÷6="),
    @("479÷6=", "653÷9="),
    @("858÷2=", "571÷4="),
    @("627÷6=", "682÷2="),
    @("540÷7=", "561÷2="),
    @("445÷9=", "600÷8="),
    @("173÷4=", "867÷7="),
    @("672÷4=", "299÷5="),
    @("108÷9=", "877÷5="),
    @("210÷5=", "953÷9="),
    @("511÷7=", "698÷7="),
    @("550÷5=", "866÷6="),
    @("997÷9=", "289÷4="),
    @("455÷9=", "439÷2="),
    @("720÷4=", "819÷5="),
    @("886÷4=", "547÷4="),
    @("220÷3=", "802÷5="),
    @("534÷7=", "538÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
